$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I3").Value = -0.4402146986448091
$ws.Range("J3").Value = 0.2205116493645886
$ws.Range("K3").Value = -0.3873567498259738
$ws.Range("L3").Value = 2.794096978151844
